# Updates market-price derived columns (H-N) across the Balmung_Profits workbook.
# Values sourced from a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow | Beeswax
$ws.Range("H12").Value = 4630159
$ws.Range("I12").Value = 5208691.5
$ws.Range("J12").Value = 1900
$ws.Range("K12").Value = 5208691.5
$ws.Range("L12").Value = 1900
$ws.Range("M12").Value = -5208521.5
$ws.Range("N12").Value = -2240

# Row 92: Whinier than the Sword | Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1563340.5
$ws.Range("I92").Value = 822801
$ws.Range("K92").Value = 822801
$ws.Range("M92").Value = -821553

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 2958.3684
$ws.Range("I98").Value = 2954.4375
$ws.Range("J98").Value = 2979.3333
$ws.Range("K98").Value = 2954.4375
$ws.Range("L98").Value = 2979.3333
$ws.Range("M98").Value = -1456.4375
$ws.Range("N98").Value = -5975.3333

# Row 107: Another Man's Ink | Enchanted Truegold Ink
$ws.Range("H107").Value = 1666.7778
$ws.Range("I107").Value = 1438.0625
$ws.Range("K107").Value = 1438.0625
$ws.Range("M107").Value = 481.9375

# Row 111: An Eye for Healing | Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 3061.5386
$ws.Range("I111").Value = 2739.125
$ws.Range("K111").Value = 8217.375
$ws.Range("M111").Value = -5150.375

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 2958.3684
$ws.Range("I122").Value = 2954.4375
$ws.Range("J122").Value = 2979.3333
$ws.Range("K122").Value = 8863.3125
$ws.Range("L122").Value = 8937.999899999999
$ws.Range("M122").Value = -6413.3125
$ws.Range("N122").Value = -13837.9999

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Range("H133").Value = 130790.836
$ws.Range("J133").Value = 130790.836
$ws.Range("L133").Value = 130790.836
$ws.Range("N133").Value = -140910.836

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 3053.875
$ws.Range("I135").Value = 3405.1667
$ws.Range("K135").Value = 30646.5003
$ws.Range("M135").Value = -28111.5003

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 5557570
$ws.Range("I137").Value = 1116.3334
$ws.Range("K137").Value = 3349.0002
$ws.Range("M137").Value = -799.0001999999999

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1082.125
$ws.Range("I141").Value = 951
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 2853
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 2327
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 47309.87
$ws.Range("I45").Value = 56085.895
$ws.Range("K45").Value = 56085.895
$ws.Range("M45").Value = -55708.895

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 3275434.2
$ws.Range("I61").Value = 8115.6665
$ws.Range("K61").Value = 8115.6665
$ws.Range("M61").Value = -7903.6665

# Row 64: Don't Scuttle with Scuta | Mythrite Scutum
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67: Shielded by Bureaucracy (L) | Mythrite Scutum
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 590167.3
$ws.Range("I74").Value = 4189.8945
$ws.Range("J74").Value = 1176144.8
$ws.Range("K74").Value = 4189.8945
$ws.Range("L74").Value = 1176144.8
$ws.Range("M74").Value = -3315.8945
$ws.Range("N74").Value = -1177892.8

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 590167.3
$ws.Range("I77").Value = 4189.8945
$ws.Range("J77").Value = 1176144.8
$ws.Range("K77").Value = 20949.4725
$ws.Range("L77").Value = 5880724
$ws.Range("M77").Value = -16581.4725
$ws.Range("N77").Value = -5889460

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 5723.522
$ws.Range("I97").Value = 6201.85
$ws.Range("K97").Value = 6201.85
$ws.Range("M97").Value = -5705.85

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 1739.0513
$ws.Range("I132").Value = 1245.1724
$ws.Range("K132").Value = 3735.5172
$ws.Range("M132").Value = -1205.5172

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3275434.2
$ws.Range("I136").Value = 8115.6665
$ws.Range("K136").Value = 24346.9995
$ws.Range("M136").Value = -21796.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 8883.272000000001
$ws.Range("I107").Value = 10525.64
$ws.Range("K107").Value = 10525.64
$ws.Range("M107").Value = -8605.639999999999

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 25716616
$ws.Range("I134").Value = 1811.8276
$ws.Range("K134").Value = 5435.4828
$ws.Range("M134").Value = -2900.4828

$ws = $wb.Worksheets.Item("CRP")
# Row 3: Touch and Heal | Maple Pattens
$ws.Range("H3").Value = 6969690
$ws.Range("I3").Value = 6969690
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6969690
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -6969577
$ws.Range("N3").ClearContents()

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 3185
$ws.Range("I31").Value = 1576.5
$ws.Range("K31").Value = 1576.5
$ws.Range("M31").Value = -1281.5

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 3185
$ws.Range("I34").Value = 1576.5
$ws.Range("K34").Value = 1576.5
$ws.Range("M34").Value = -1374.5

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2486.1667
$ws.Range("I58").Value = 1902.8572
$ws.Range("J58").Value = 3302.8
$ws.Range("K58").Value = 1902.8572
$ws.Range("L58").Value = 3302.8
$ws.Range("M58").Value = -1699.8572
$ws.Range("N58").Value = -3708.8

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 15866.357
$ws.Range("I86").Value = 10295.444
$ws.Range("J86").Value = 25894
$ws.Range("K86").Value = 10295.444
$ws.Range("L86").Value = 25894
$ws.Range("M86").Value = -9172.444
$ws.Range("N86").Value = -28140

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 15866.357
$ws.Range("I89").Value = 10295.444
$ws.Range("J89").Value = 25894
$ws.Range("K89").Value = 51477.22
$ws.Range("L89").Value = 129470
$ws.Range("M89").Value = -45861.22
$ws.Range("N89").Value = -140702

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 2221
$ws.Range("I105").Value = 953.6
$ws.Range("J105").Value = 4333.3335
$ws.Range("K105").Value = 953.6
$ws.Range("L105").Value = 4333.3335
$ws.Range("M105").Value = 793.4
$ws.Range("N105").Value = -7827.3335

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 3640.394
$ws.Range("I122").Value = 2310.1738
$ws.Range("K122").Value = 6930.5214
$ws.Range("M122").Value = -4480.5214

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2486.1667
$ws.Range("I136").Value = 1902.8572
$ws.Range("J136").Value = 3302.8
$ws.Range("K136").Value = 5708.571599999999
$ws.Range("L136").Value = 9908.400000000001
$ws.Range("M136").Value = -3158.571599999999
$ws.Range("N136").Value = -15008.4

$ws = $wb.Worksheets.Item("CUL")
# Row 44: No More Dumpster Diving | Knight's Bread
$ws.Range("H44").Value = 3374.25
$ws.Range("I44").Value = 1948.5
$ws.Range("J44").Value = 4800
$ws.Range("K44").Value = 5845.5
$ws.Range("L44").Value = 14400
$ws.Range("M44").Value = -5447.5
$ws.Range("N44").Value = -15196

# Row 127: A Stickler for Carrots | Carrot Nibbles
$ws.Range("H127").Value = 35499.5
$ws.Range("J127").Value = 35499.5
$ws.Range("L127").Value = 106498.5
$ws.Range("N127").Value = -116418.5

# Row 133: Friends Are Food | Boiled Alpaca Steak
$ws.Range("H133").Value = 4074.1333
$ws.Range("I133").Value = 3842.6667
$ws.Range("K133").Value = 11528.0001
$ws.Range("M133").Value = -6468.000100000001

# Row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value = 1328
$ws.Range("I136").Value = 1328
$ws.Range("K136").Value = 3984
$ws.Range("M136").Value = 1116

# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 4474.067
$ws.Range("J137").Value = 8898.799999999999
$ws.Range("L137").Value = 26696.4
$ws.Range("N137").Value = -36896.39999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 3: Needful Rings | Copper Wristlets
$ws.Range("H3").Value = 508400
$ws.Range("I3").Value = 760000
$ws.Range("K3").Value = 760000
$ws.Range("M3").Value = -759884

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 10211323
$ws.Range("I132").Value = 1380.9333
$ws.Range("K132").Value = 4142.7999
$ws.Range("M132").Value = -1612.7999

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 9094786
$ws.Range("I61").Value = 11768606
$ws.Range("K61").Value = 11768606
$ws.Range("M61").Value = -11768404

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 9094786
$ws.Range("I113").Value = 11768606
$ws.Range("K113").Value = 11768606
$ws.Range("M113").Value = -11766436

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1734.4375
$ws.Range("J81").Value = 1587
$ws.Range("L81").Value = 3174
$ws.Range("N81").Value = -5296

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1734.4375
$ws.Range("J84").Value = 1587
$ws.Range("L84").Value = 15870
$ws.Range("N84").Value = -26478

# Row 95: Duress Rehearsal | Ruby Cotton Fingerless Gloves of Casting
$ws.Range("H95").Value = 65944
$ws.Range("J95").Value = 65944
$ws.Range("L95").Value = 65944
$ws.Range("N95").Value = -71436

# Row 100: Of Great Import | Kudzu Thread
$ws.Range("H100").Value = 948.5
$ws.Range("I100").Value = 989.1
$ws.Range("J100").Value = 745.5
$ws.Range("K100").Value = 1978.2
$ws.Range("L100").Value = 1491
$ws.Range("M100").Value = -1437.2
$ws.Range("N100").Value = -2573

# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 52632010
$ws.Range("I107").Value = 528
$ws.Range("K107").Value = 1584
$ws.Range("M107").Value = 336

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2573.5
$ws.Range("I122").Value = 2836.6155
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 8509.8465
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -6059.8465
$ws.Range("N122").Value = -9200.0002
